# Removed Test Case Inter-Dependency
#
# - Product name (A1/B1 on both sheets) is changed to a new, distinct value
#   so the test case does not depend on another run's product name.
# - Short name (A2/B2 on ProductLoanInput) switches from the numeric 4211
#   to a distinct text value "421w".
# - The active/selected sheet moves from "ProductLoanInput" to
#   "ProductLoanOutput".

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "4211-RBI-EI-DB-DL-REC-RNI-FEE+INT-FFConMONTHLYonDAY25-FIFC-1-FFROP-DAILY-FIFR-1-MD-TR-1-ON-PER-1st"
$newShortName   = "421w"

# Update product name on both sheets
$wsInput.Range("B1").Value  = $newProductName
$wsOutput.Range("B1").Value = $newProductName

# Update short name on the input sheet
$wsInput.Range("B2").Value = $newShortName

# Switch the active sheet / selection so ProductLoanOutput is the one shown
# when the workbook is reopened.
$wsInput.Range("A30").Select() | Out-Null
$wsOutput.Activate() | Out-Null
$wsOutput.Range("B1").Select() | Out-Null
